# Clean up footnote markers (e.g. " [5]") and embedded line breaks in
# vaccine/brand-name labels across every worksheet of the workbook.
#
# Rule applied to every text cell:
#   1. Strip any "[<digits>]" footnote-reference token.
#   2. Collapse embedded newlines (two-line wrapped labels) into a single
#      space so the label reads as one line.

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $ur = $ws.UsedRange
    $rowCount = $ur.Rows.Count
    $colCount = $ur.Columns.Count

    $firstRow = $ur.Row
    $firstCol = $ur.Column

    for ($r = 0; $r -lt $rowCount; $r++) {
        for ($c = 0; $c -lt $colCount; $c++) {
            $cell = $ws.Cells.Item($firstRow + $r, $firstCol + $c)
            $val = $cell.Value2

            if ($val -ne $null -and $val -is [string]) {
                $newVal = $val -replace '\[\d+\]', ''
                $newVal = $newVal -replace "`r`n", ' '
                $newVal = $newVal -replace "`n", ' '

                if ($newVal -ne $val) {
                    $cell.Value = $newVal
                }
            }
        }
    }
}
